# Apply gameflow.xlsx changes: add mid-battle interrupt/cutscene rows
# (health-threshold triggers for dialogue interrupts during battle),
# per "Starting implementation of mid-battle cutscenes." commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4
$ws.Range("A2").Value = 'START_SCENE'
$ws.Range("B2").Value = 'INTRO'
$ws.Range("A3").Value = 'END_SCENE'
$ws.Range("A4").Value = 'START_SCENE'
$ws.Range("B4").Value = 'CUTSCENE'
$ws.Range("C4").Value = 'jazzy_retro_battle_theme'
$ws.Range("A5").Value = 'NPC'
$ws.Range("B5").Value = 'Tanuki'
$ws.Range("A6").Value = 'NPC'
$ws.Range("B6").Value = 'Frog'
$ws.Range("A7").Value = 'DIALOGUE'
$ws.Range("B7").Value = 'Tanuki'
$ws.Range("C7").Value = 'Hey you! You''re walking in the wrong part of town.'
$ws.Range("D7").Value = 'tanuki_mario'
$ws.Range("E7").Value = 'jazzy_retro_battle_theme'
$ws.Range("A8").Value = 'DIALOGUE'
$ws.Range("B8").Value = 'Frog'
$ws.Range("C8").Value = 'Ribbit Ribbit! (Yeah frog-face! Wrong part of town!)'
$ws.Range("D8").Value = 'frog_mario'
$ws.Range("E8").Value = 'frogs'
$ws.Range("A9").Value = 'DIALOGUE'
$ws.Range("B9").Value = 'Tanuki'
$ws.Range("C9").Value = 'Let''s get em!'
$ws.Range("D9").Value = 'tanuki_mario'
$ws.Range("E9").Value = 'jazzy_retro_battle_theme'
$ws.Range("A10").Value = 'END_SCENE'
$ws.Range("A11").Value = 'START_SCENE'
$ws.Range("B11").Value = 'BATTLE'
$ws.Range("A12").Value = 'MUSIC'
$ws.Range("B12").Value = 'sample_incomplete_war_remix'
$ws.Range("A13").Value = 'ENEMY'
$ws.Range("B13").Value = 'Slime'
$ws.Range("A14").Value = 'ENEMY'
$ws.Range("B14").Value = 'The Evil Eye'
$ws.Range("A15").Value = 'ENEMY'
$ws.Range("B15").Value = 'Tanuki'
$ws.Range("A16").Value = 'INTERRUPT'
$ws.Range("B16").Value = 'MIDDLE_HEALTH'
$ws.Range("C16").Value = 0.5
$ws.Range("D16").Value = 'L'
$ws.Range("A17").Value = 'DIALOGUE'
$ws.Range("B17").Value = 'Tanuki'
$ws.Range("C17").Value = 'Ribbit! (Wow, you''re actually killing the Evil Eye!)'
$ws.Range("D17").Value = 'frog_mario'
$ws.Range("E17").Value = '_'
$ws.Range("A18").Value = 'END_INTERRUPT'
$ws.Range("A19").Value = 'INTERRUPT'
$ws.Range("B19").Value = 'MIDDLE_HEALTH'
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 'R'
$ws.Range("A20").Value = 'DIALOGUE'
$ws.Range("B20").Value = 'Tanuki'
$ws.Range("C20").Value = 'Did you just kill the Evil Eye? Did- Did you win?'
$ws.Range("D20").Value = 'tanuki_mario'
$ws.Range("E20").Value = '_'
$ws.Range("A21").Value = 'DIALOGUE'
$ws.Range("B21").Value = 'Tanuki'
$ws.Range("C21").Value = 'Holy crap.'
$ws.Range("D21").Value = 'tanuki_mario'
$ws.Range("E21").Value = '_'
$ws.Range("A22").Value = 'INTERRUPT_END'
$ws.Range("A23").Value = 'INTERRUPT'
$ws.Range("B23").Value = 'RIGHT_HEALTH'
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 'R'
$ws.Range("A24").Value = 'DIALOGUE'
$ws.Range("B24").Value = 'Tanuki'
$ws.Range("C24").Value = 'Ah! You have defeated me!'
$ws.Range("D24").Value = 'tanuki_mario'
$ws.Range("E24").Value = '_'
$ws.Range("A25").Value = 'END_INTERRUPT'
$ws.Range("A26").Value = 'INTERRUPT'
$ws.Range("B26").Value = 'PLAYER_HEALTH'
$ws.Range("C26").Value = 0.75
$ws.Range("D26").Value = 'LR'
$ws.Range("A27").Value = 'DIALOGUE'
$ws.Range("B27").Value = 'Tanuki'
$ws.Range("C27").Value = 'Ha! Feel the wrath of Tanuki and Frog!'
$ws.Range("D27").Value = 'tanuki_mario'
$ws.Range("E27").Value = '_'
$ws.Range("A28").Value = 'DIALOGUE'
$ws.Range("B28").Value = 'Frog'
$ws.Range("C28").Value = 'Ribbit Ribbit Rrrrrrribbit! (Frog is actually just making frog sounds right now)'
$ws.Range("D28").Value = 'frog_mario'
$ws.Range("E28").Value = '_'
$ws.Range("A29").Value = 'END_INTERRUPT'
$ws.Range("A30").Value = 'END_SCENE'
$ws.Range("A31").Value = 'START_SCENE'
$ws.Range("B31").Value = 'CUTSCENE'
$ws.Range("A32").Value = 'NPC'
$ws.Range("B32").Value = 'Tanuki'
$ws.Range("A33").Value = 'NPC'
$ws.Range("B33").Value = 'Frog'
$ws.Range("A34").Value = 'DIALOGUE'
$ws.Range("B34").Value = 'Tanuki'
$ws.Range("C34").Value = 'Ugh! How are you so powerful? You even killed Frog!'
$ws.Range("D34").Value = 'tanuki_mario'
$ws.Range("E34").Value = 'second_hand'
$ws.Range("A35").Value = 'DIALOGUE'
$ws.Range("B35").Value = 'Frog'
$ws.Range("C35").Value = '(He''s not moving)'
$ws.Range("D35").Value = 'frog_mario'
$ws.Range("E35").Value = '_'
$ws.Range("A36").Value = 'DIALOGUE'
$ws.Range("B36").Value = 'Tanuki'
$ws.Range("C36").Value = 'Bleh. (ded)'
$ws.Range("D36").Value = 'tanuki_mario'
$ws.Range("E36").Value = '_'
$ws.Range("A37").Value = 'END_SCENE'
$ws.Range("A38").Value = 'END_GAME'

# Update the view state to match the edited region (best effort).
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("B14").Select()
